$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; this shifts the existing rows 15-67 down to 16-68,
# preserving their data (matching the target diff's row renumbering).
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new weekly record.
$ws.Cells.Item(15, 1).Value = 1
$ws.Cells.Item(15, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(15, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(15, 4).Value = 44592
$ws.Cells.Item(15, 5).Value = 15
$ws.Cells.Item(15, 6).Value = 100112038
$ws.Cells.Item(15, 7).Value = "Cebollín baby"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 250
$ws.Cells.Item(15, 11).Value = 3000
$ws.Cells.Item(15, 12).Value = 3500
$ws.Cells.Item(15, 13).Value = 3250
$ws.Cells.Item(15, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(15, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(15, 16).Value = 1625
$ws.Cells.Item(15, 17).Value = 2
$ws.Cells.Item(15, 18).Value = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Cells.Item(15, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
